$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose text looks like a plain number must be forced to Text format
# so Excel keeps them as strings (matching the source data which stores
# prices as text), then the style is reset back to Normal so no stray
# number-format style is left attached to the cell.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9995'
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '238.64'
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4731'
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2638'
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06211'
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07057'
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.5921'
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.414'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '76.22'
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.000'
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000006809'
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.55'
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.548'
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.767'
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.321'
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '134.97'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.25'
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.404'
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '108.08'
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.007'
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.685'
$ws.Range("D31").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04431'
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.614'
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9784'
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6195'
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9334'
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '113.78'
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.412'
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.915'
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.01478'
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.324'
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3821'
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '6.292'
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05284'
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '30.39'
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.705'
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3367'
$ws.Range("D51").Style = "Normal"

# Remaining cells already parse as text in Excel (non-numeric content,
# e.g. thousand-grouped prices or percentage strings), so a plain value
# assignment is sufficient.
$ws.Range("D2").Value = '26.359.55'
$ws.Range("E2").Value = '  +2.95%  '
$ws.Range("D3").Value = '1.718.68'
$ws.Range("E3").Value = '  +3.24%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("E5").Value = '  +0.87%  '
$ws.Range("E6").Value = '  +0.00%  '
$ws.Range("E7").Value = '  -1.43%  '
$ws.Range("E8").Value = '  +0.49%  '
$ws.Range("E9").Value = '  +0.86%  '
$ws.Range("D10").Value = '1.716.50'
$ws.Range("E10").Value = '  +3.11%  '
$ws.Range("E11").Value = '  -0.43%  '
$ws.Range("E12").Value = '  +3.51%  '
$ws.Range("E13").Value = '  -0.20%  '
$ws.Range("E14").Value = '  +0.57%  '
$ws.Range("E15").Value = '  +2.46%  '
$ws.Range("E16").Value = '  +0.03%  '
$ws.Range("E17").Value = '  +0.06%  '
$ws.Range("D18").Value = '26.345.93'
$ws.Range("E18").Value = '  +2.90%  '
$ws.Range("E19").Value = '  +0.72%  '
$ws.Range("D21").Value = '1.935.83'
$ws.Range("E21").Value = '  +3.09%  '
$ws.Range("E22").Value = '  +2.43%  '
$ws.Range("E23").Value = '  +1.12%  '
$ws.Range("E24").Value = '  +0.24%  '
$ws.Range("E25").Value = '  +0.11%  '
$ws.Range("E26").Value = '  +1.26%  '
$ws.Range("E27").Value = '  +0.13%  '
$ws.Range("E28").Value = '  +3.25%  '
$ws.Range("E29").Value = '  +3.63%  '
$ws.Range("E30").Value = '  +1.34%  '
$ws.Range("E31").Value = '  +0.44%  '
$ws.Range("E32").Value = '  +1.16%  '
$ws.Range("E33").Value = '  +2.54%  '
$ws.Range("E35").Value = '  +2.87%  '
$ws.Range("E36").Value = '  +1.21%  '
$ws.Range("E37").Value = '  +9.14%  '
$ws.Range("E38").Value = '  +16.02%  '
$ws.Range("E39").Value = '  -7.51%  '
$ws.Range("E40").Value = '  +1.57%  '
$ws.Range("E41").Value = '  +0.00%  '
$ws.Range("E42").Value = '  -2.00%  '
$ws.Range("E43").Value = '  +13.01%  '
$ws.Range("E44").Value = '  +1.37%  '
$ws.Range("E45").Value = '  +4.28%  '
$ws.Range("E46").Value = '  +1.15%  '
$ws.Range("E47").Value = '  +0.41%  '
$ws.Range("E48").Value = '  +3.03%  '
$ws.Range("E49").Value = '  +4.67%  '
$ws.Range("E50").Value = '  +1.61%  '
$ws.Range("E51").Value = '  +0.79%  '
